$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "Davide Zeni"
$ws.Range("B65").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C65").Value = "Andrea  Roveda  | Pinguini Trentini"
$ws.Range("D65").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E65").Value = "Luca Perenzoni | CGB Gamberoni"
$ws.Range("F65").Value = "Jacopo  Chemini | IMONTAGNA"
